# Update the "From" value of the R30 rule row (row 10) on the Rules sheet
# from 18 to 1, as described by the diff (cell C10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
